$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.404689
$ws.Range("H2").Value = 16.214067
$ws.Range("I2").Value = 0.08747859311663772
$ws.Range("J2").Value = 0.09021076583983562
$ws.Range("M2").Value = 2.906846333333333
$ws.Range("N2").Value = 8.720538999999999
$ws.Range("O2").Value = 0.005520525738044089
$ws.Range("P2").Value = 0.005624540846623205
$ws.Range("Q2").Value = 15.710600402457
$ws.Range("R2").Value = 141.395403622113
$ws.Range("S2").Value = 0.000482927824828285
$ws.Range("T2").Value = 0.0005073941372713168
$ws.Range("G3").Value = 5.404689
$ws.Range("H3").Value = 16.214067
$ws.Range("I3").Value = 0.08747859311663772
$ws.Range("J3").Value = 0.09021076583983562
$ws.Range("O3").Value = 0.3528665483720876
$ws.Range("P3").Value = 0.3595150912979765
$ws.Range("Q3").Value = 1004.206048468179
$ws.Range("R3").Value = 9037.854436213611
$ws.Range("S3").Value = 0.03086826920951422
$ws.Range("T3").Value = 0.03243213171696888
$ws.Range("G4").Value = 5.404689
$ws.Range("H4").Value = 16.214067
$ws.Range("I4").Value = 0.08747859311663772
$ws.Range("J4").Value = 0.09021076583983562
$ws.Range("M4").Value = 137.0717086666666
$ws.Range("N4").Value = 411.2151259999999
$ws.Range("O4").Value = 0.2603191943704447
$ws.Range("P4").Value = 0.2652240042658267
$ws.Range("Q4").Value = 740.8299560419379
$ws.Range("R4").Value = 6667.469604377441
$ws.Range("S4").Value = 0.02277235688478306
$ws.Range("T4").Value = 0.02392606054392805
$ws.Range("G5").Value = 5.404689
$ws.Range("H5").Value = 16.214067
$ws.Range("I5").Value = 0.08747859311663772
$ws.Range("J5").Value = 0.09021076583983562
$ws.Range("M5").Value = 29.2127365
$ws.Range("N5").Value = 58.425473
$ws.Range("O5").Value = 0.05547925319534149
$ws.Range("P5").Value = 0.03768304451958546
$ws.Range("Q5").Value = 157.8857556214485
$ws.Range("R5").Value = 947.3145337286909
$ws.Range("S5").Value = 0.004853247016690201
$ws.Range("T5").Value = 0.003399416305288424
$ws.Range("G6").Value = 5.404689
$ws.Range("H6").Value = 16.214067
$ws.Range("I6").Value = 0.08747859311663772
$ws.Range("J6").Value = 0.09021076583983562
$ws.Range("M6").Value = 171.5584106666666
$ws.Range("N6").Value = 514.6752319999999
$ws.Range("O6").Value = 0.3258144783240821
$ws.Range("P6").Value = 0.331953319069988
$ws.Range("Q6").Value = 927.2198549876159
$ws.Range("R6").Value = 8344.978694888543
$ws.Range("S6").Value = 0.02850179218082196
$ws.Range("T6").Value = 0.02994576313637893
$ws.Range("I7").Value = 0.2448858138641327
$ws.Range("J7").Value = 0.2525342032254661
$ws.Range("M7").Value = 2.906846333333333
$ws.Range("N7").Value = 8.720538999999999
$ws.Range("O7").Value = 0.005520525738044089
$ws.Range("P7").Value = 0.005624540846623205
$ws.Range("Q7").Value = 43.97993873449855
$ws.Range("R7").Value = 395.819448610487
$ws.Range("S7").Value = 0.001351898438318819
$ws.Range("T7").Value = 0.00142038894121108
$ws.Range("I8").Value = 0.2448858138641327
$ws.Range("J8").Value = 0.2525342032254661
$ws.Range("O8").Value = 0.3528665483720876
$ws.Range("P8").Value = 0.3595150912979765
$ws.Range("S8").Value = 0.08641201188352604
$ws.Range("T8").Value = 0.09078985712846521
$ws.Range("I9").Value = 0.2448858138641327
$ws.Range("J9").Value = 0.2525342032254661
$ws.Range("M9").Value = 137.0717086666666
$ws.Range("N9").Value = 411.2151259999999
$ws.Range("O9").Value = 0.2603191943704447
$ws.Range("P9").Value = 0.2652240042658267
$ws.Range("Q9").Value = 2073.864476516773
$ws.Range("R9").Value = 18664.78028865096
$ws.Range("S9").Value = 0.0637484777778617
$ws.Range("T9").Value = 0.06697813259353817
$ws.Range("I10").Value = 0.2448858138641327
$ws.Range("J10").Value = 0.2525342032254661
$ws.Range("M10").Value = 29.2127365
$ws.Range("N10").Value = 58.425473
$ws.Range("O10").Value = 0.05547925319534149
$ws.Range("P10").Value = 0.03768304451958546
$ws.Range("Q10").Value = 441.9822082799182
$ws.Range("R10").Value = 2651.893249679509
$ws.Range("S10").Value = 0.01358608207131549
$ws.Range("T10").Value = 0.00951625762286328
$ws.Range("I11").Value = 0.2448858138641327
$ws.Range("J11").Value = 0.2525342032254661
$ws.Range("M11").Value = 171.5584106666666
$ws.Range("N11").Value = 514.6752319999999
$ws.Range("O11").Value = 0.3258144783240821
$ws.Range("P11").Value = 0.331953319069988
$ws.Range("Q11").Value = 2595.640610233362
$ws.Range("R11").Value = 23360.76549210025
$ws.Range("S11").Value = 0.07978734369311069
$ws.Range("T11").Value = 0.08382956693938835
$ws.Range("G12").Value = 17.564497
$ws.Range("H12").Value = 52.69349099999999
$ws.Range("I12").Value = 0.2842934138044583
$ws.Range("J12").Value = 0.2931725999334087
$ws.Range("M12").Value = 2.906846333333333
$ws.Range("N12").Value = 8.720538999999999
$ws.Range("O12").Value = 0.005520525738044089
$ws.Range("P12").Value = 0.005624540846623205
$ws.Range("Q12").Value = 51.05729370129432
$ws.Range("R12").Value = 459.5156433116489
$ws.Range("S12").Value = 0.001569449108063931
$ws.Range("T12").Value = 0.001648961263436181
$ws.Range("G13").Value = 17.564497
$ws.Range("H13").Value = 52.69349099999999
$ws.Range("I13").Value = 0.2842934138044583
$ws.Range("J13").Value = 0.2931725999334087
$ws.Range("O13").Value = 0.3528665483720876
$ws.Range("P13").Value = 0.3595150912979765
$ws.Range("Q13").Value = 3263.5317454346
$ws.Range("R13").Value = 29371.7857089114
$ws.Range("S13").Value = 0.1003176356540968
$ws.Range("T13").Value = 0.1053999740311246
$ws.Range("G14").Value = 17.564497
$ws.Range("H14").Value = 52.69349099999999
$ws.Range("I14").Value = 0.2842934138044583
$ws.Range("J14").Value = 0.2931725999334087
$ws.Range("M14").Value = 137.0717086666666
$ws.Range("N14").Value = 411.2151259999999
$ws.Range("O14").Value = 0.2603191943704447
$ws.Range("P14").Value = 0.2652240042658267
$ws.Range("Q14").Value = 2407.59561566054
$ws.Range("R14").Value = 21668.36054094486
$ws.Range("S14").Value = 0.07400703244640003
$ws.Range("T14").Value = 0.07775641089536189
$ws.Range("G15").Value = 17.564497
$ws.Range("H15").Value = 52.69349099999999
$ws.Range("I15").Value = 0.2842934138044583
$ws.Range("J15").Value = 0.2931725999334087
$ws.Range("M15").Value = 29.2127365
$ws.Range("N15").Value = 58.425473
$ws.Range("O15").Value = 0.05547925319534149
$ws.Range("P15").Value = 0.03768304451958546
$ws.Range("Q15").Value = 513.1070226160405
$ws.Range("R15").Value = 3078.642135696242
$ws.Range("S15").Value = 0.01577238628622553
$ws.Range("T15").Value = 0.01104763613521326
$ws.Range("G16").Value = 17.564497
$ws.Range("H16").Value = 52.69349099999999
$ws.Range("I16").Value = 0.2842934138044583
$ws.Range("J16").Value = 0.2931725999334087
$ws.Range("M16").Value = 171.5584106666666
$ws.Range("N16").Value = 514.6752319999999
$ws.Range("O16").Value = 0.3258144783240821
$ws.Range("P16").Value = 0.331953319069988
$ws.Range("Q16").Value = 3013.337189479434
$ws.Range("R16").Value = 27120.03470531491
$ws.Range("S16").Value = 0.09262691030967198
$ws.Range("T16").Value = 0.09731961760827278
$ws.Range("G17").Value = 5.613580499999999
$ws.Range("H17").Value = 11.227161
$ws.Range("I17").Value = 0.09085964511315853
$ws.Range("J17").Value = 0.06246494429911598
$ws.Range("M17").Value = 2.906846333333333
$ws.Range("N17").Value = 8.720538999999999
$ws.Range("O17").Value = 0.005520525738044089
$ws.Range("P17").Value = 0.005624540846623205
$ws.Range("Q17").Value = 16.31781589329649
$ws.Range("R17").Value = 97.90689535977897
$ws.Range("S17").Value = 0.0005015930093967434
$ws.Range("T17").Value = 0.0003513366306924212
$ws.Range("G18").Value = 5.613580499999999
$ws.Range("H18").Value = 11.227161
$ws.Range("I18").Value = 0.09085964511315853
$ws.Range("J18").Value = 0.06246494429911598
$ws.Range("O18").Value = 0.3528665483720876
$ws.Range("P18").Value = 0.3595150912979765
$ws.Range("Q18").Value = 1043.018662436085
$ws.Range("R18").Value = 6258.111974616512
$ws.Range("S18").Value = 0.03206132935739307
$ws.Range("T18").Value = 0.0224570901526197
$ws.Range("G19").Value = 5.613580499999999
$ws.Range("H19").Value = 11.227161
$ws.Range("I19").Value = 0.09085964511315853
$ws.Range("J19").Value = 0.06246494429911598
$ws.Range("M19").Value = 137.0717086666666
$ws.Range("N19").Value = 411.2151259999999
$ws.Range("O19").Value = 0.2603191943704447
$ws.Range("P19").Value = 0.2652240042658267
$ws.Range("Q19").Value = 769.4630708728807
$ws.Range("R19").Value = 4616.778425237285
$ws.Range("S19").Value = 0.02365250961664194
$ws.Range("T19").Value = 0.01656720265325336
$ws.Range("G20").Value = 5.613580499999999
$ws.Range("H20").Value = 11.227161
$ws.Range("I20").Value = 0.09085964511315853
$ws.Range("J20").Value = 0.06246494429911598
$ws.Range("M20").Value = 29.2127365
$ws.Range("N20").Value = 58.425473
$ws.Range("O20").Value = 0.05547925319534149
$ws.Range("P20").Value = 0.03768304451958546
$ws.Range("Q20").Value = 163.9880479680382
$ws.Range("R20").Value = 655.9521918721529
$ws.Range("S20").Value = 0.005040825256471794
$ws.Range("T20").Value = 0.002353869276937013
$ws.Range("G21").Value = 5.613580499999999
$ws.Range("H21").Value = 11.227161
$ws.Range("I21").Value = 0.09085964511315853
$ws.Range("J21").Value = 0.06246494429911598
$ws.Range("M21").Value = 171.5584106666666
$ws.Range("N21").Value = 514.6752319999999
$ws.Range("O21").Value = 0.3258144783240821
$ws.Range("P21").Value = 0.331953319069988
$ws.Range("Q21").Value = 963.0569487293918
$ws.Range("R21").Value = 5778.341692376351
$ws.Range("S21").Value = 0.02960338787325498
$ws.Range("T21").Value = 0.02073544558561348
$ws.Range("G22").Value = 18.07044533333334
$ws.Range("H22").Value = 54.211336
$ws.Range("I22").Value = 0.2924825341016128
$ws.Range("J22").Value = 0.3016174867021735
$ws.Range("M22").Value = 2.906846333333333
$ws.Range("N22").Value = 8.720538999999999
$ws.Range("O22").Value = 0.005520525738044089
$ws.Range("P22").Value = 0.005624540846623205
$ws.Range("Q22").Value = 52.52800775890044
$ws.Range("R22").Value = 472.7520698301039
$ws.Range("S22").Value = 0.001614657357436311
$ws.Range("T22").Value = 0.001696459874012207
$ws.Range("G23").Value = 18.07044533333334
$ws.Range("H23").Value = 54.211336
$ws.Range("I23").Value = 0.2924825341016128
$ws.Range("J23").Value = 0.3016174867021735
$ws.Range("O23").Value = 0.3528665483720876
$ws.Range("P23").Value = 0.3595150912979765
$ws.Range("Q23").Value = 3357.53833425881
$ws.Range("R23").Value = 30217.84500832929
$ws.Range("S23").Value = 0.1032073022675575
$ws.Range("T23").Value = 0.1084360382687981
$ws.Range("G24").Value = 18.07044533333334
$ws.Range("H24").Value = 54.211336
$ws.Range("I24").Value = 0.2924825341016128
$ws.Range("J24").Value = 0.3016174867021735
$ws.Range("M24").Value = 137.0717086666666
$ws.Range("N24").Value = 411.2151259999999
$ws.Range("O24").Value = 0.2603191943704447
$ws.Range("P24").Value = 0.2652240042658267
$ws.Range("Q24").Value = 2476.946818207593
$ws.Range("R24").Value = 22292.52136386833
$ws.Range("S24").Value = 0.07613881764475794
$ws.Range("T24").Value = 0.0799961975797452
$ws.Range("G25").Value = 18.07044533333334
$ws.Range("H25").Value = 54.211336
$ws.Range("I25").Value = 0.2924825341016128
$ws.Range("J25").Value = 0.3016174867021735
$ws.Range("M25").Value = 29.2127365
$ws.Range("N25").Value = 58.425473
$ws.Range("O25").Value = 0.05547925319534149
$ws.Range("P25").Value = 0.03768304451958546
$ws.Range("Q25").Value = 527.8871579603214
$ws.Range("R25").Value = 3167.322947761928
$ws.Range("S25").Value = 0.01622671256463848
$ws.Range("T25").Value = 0.01136586517928348
$ws.Range("G26").Value = 18.07044533333334
$ws.Range("H26").Value = 54.211336
$ws.Range("I26").Value = 0.2924825341016128
$ws.Range("J26").Value = 0.3016174867021735
$ws.Range("M26").Value = 171.5584106666666
$ws.Range("N26").Value = 514.6752319999999
$ws.Range("O26").Value = 0.3258144783240821
$ws.Range("P26").Value = 0.331953319069988
$ws.Range("Q26").Value = 3100.13688142555
$ws.Range("R26").Value = 27901.23193282995
$ws.Range("S26").Value = 0.09529504426722253
$ws.Range("T26").Value = 0.1001229258003345

Write-Host "Applied 278 cell updates"